$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.800.84"
$ws.Range("E2").Value = "'  -3.27%  "
$ws.Range("D3").Value = "'2.556.87"
$ws.Range("E3").Value = "'  -1.67%  "
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'504.58"
$ws.Range("E5").Value = "'  -3.60%  "
$ws.Range("D6").Value = "'141.42"
$ws.Range("E6").Value = "'  -8.09%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("D8").Value = "'0.552"
$ws.Range("E8").Value = "'  -5.84%  "
$ws.Range("D9").Value = "'2.560.23"
$ws.Range("E9").Value = "'  -1.73%  "
$ws.Range("D10").Value = "'6.18"
$ws.Range("E10").Value = "'  -7.59%  "
$ws.Range("E11").Value = "'  -4.19%  "
$ws.Range("E13").Value = "'  -1.03%  "
$ws.Range("D14").Value = "'3.006.26"
$ws.Range("E14").Value = "'  -1.59%  "
$ws.Range("D15").Value = "'58.815.58"
$ws.Range("E15").Value = "'  -3.28%  "
$ws.Range("D16").Value = "'20.46"
$ws.Range("E16").Value = "'  -5.14%  "
$ws.Range("E17").Value = "'  -4.79%  "
$ws.Range("D18").Value = "'2.578.27"
$ws.Range("E18").Value = "'  -0.98%  "
$ws.Range("E19").Value = "'  -5.38%  "
$ws.Range("D20").Value = "'331.44"
$ws.Range("E20").Value = "'  -6.58%  "
$ws.Range("D21").Value = "'10.01"
$ws.Range("E21").Value = "'  -5.04%  "
$ws.Range("E22").Value = "'  +0.02%  "
$ws.Range("E23").Value = "'  -4.25%  "
$ws.Range("D24").Value = "'59.45"
$ws.Range("E24").Value = "'  -2.66%  "
$ws.Range("D25").Value = "'0.403"
$ws.Range("E25").Value = "'  -5.06%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "'  +0.15%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = "'  -4.06%  "
$ws.Range("E28").Value = "'  -7.88%  "
$ws.Range("E29").Value = "'  -7.13%  "
$ws.Range("E30").Value = "'  +0.00%  "
$ws.Range("D31").Value = "'148.66"
$ws.Range("E31").Value = "'  -1.08%  "
$ws.Range("E32").Value = "'  -4.70%  "
$ws.Range("E33").Value = "'  -4.17%  "
$ws.Range("D34").Value = "'5.76"
$ws.Range("E34").Value = "'  -8.04%  "
$ws.Range("D35").Value = "'3.85"
$ws.Range("E35").Value = "'  -7.75%  "
$ws.Range("D36").Value = "'0.877"
$ws.Range("E36").Value = "'  -4.63%  "
$ws.Range("E37").Value = "'  -8.13%  "
$ws.Range("D38").Value = "'35.78"
$ws.Range("E38").Value = "'  -1.51%  "
$ws.Range("D39").Value = "'0.818"
$ws.Range("E39").Value = "'  -9.60%  "
$ws.Range("D40").Value = "'286.05"
$ws.Range("E40").Value = "'  -3.73%  "
$ws.Range("E42").Value = "'  -7.49%  "
$ws.Range("E43").Value = "'  +0.02%  "
$ws.Range("D44").Value = "'0.0979"
$ws.Range("E44").Value = "'  -3.30%  "
$ws.Range("D45").Value = "'0.605"
$ws.Range("E45").Value = "'  -2.97%  "
$ws.Range("D46").Value = "'0.0527"
$ws.Range("E46").Value = "'  -5.52%  "
$ws.Range("E47").Value = "'  +0.03%  "
$ws.Range("D48").Value = "'18.56"
$ws.Range("E48").Value = "'  -4.97%  "
$ws.Range("E49").Value = "'  -5.39%  "
$ws.Range("D50").Value = "'4.50"
$ws.Range("E50").Value = "'  -8.17%  "
$ws.Range("D51").Value = "'1.886.00"
$ws.Range("E51").Value = "'  -4.00%  "
